$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 and Row 4 swap their distinguishing field values (A, K, Q, R, Z, AB, AC).
$cols = @("A", "K", "Q", "R", "Z", "AB", "AC")

foreach ($col in $cols) {
    $addr3 = "$col" + "3"
    $addr4 = "$col" + "4"
    $v3 = $ws.Range($addr3).Value()
    $v4 = $ws.Range($addr4).Value()
    $ws.Range($addr3).Value = $v4
    $ws.Range($addr4).Value = $v3
}
